$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 4956
$ws.Range("F5").Value = 2846
$ws.Range("F7").Value = 1961
$ws.Range("F9").Value = 1753
$ws.Range("F10").Value = 763
$ws.Range("F13").Value = 444
$ws.Range("F15").Value = 307
$ws.Range("F16").Value = 13
$ws.Range("F18").Value = 82
$ws.Range("F19").Value = 1056
$ws.Range("F22").Value = 691
$ws.Range("F23").Value = 765
$ws.Range("F25").Value = 17
$ws.Range("F27").Value = 586
$ws.Range("F28").Value = 72
$ws.Range("F29").Value = 1670
$ws.Range("F30").Value = 1747
$ws.Range("F31").Value = 431
$ws.Range("F33").Value = 1629
$ws.Range("F34").Value = 236
$ws.Range("F35").Value = 2440
$ws.Range("F36").Value = 429
$ws.Range("F38").Value = 637
$ws.Range("F40").Value = 75
$ws.Range("F43").Value = 1536
$ws.Range("F44").Value = 244
$ws.Range("F49").Value = 126

# --- Sheet "演出" (Performance) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 7
$ws.Range("F4").Value = 111
$ws.Range("F12").Value = 53

# --- Sheet "全部类型" (All types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 4956
$ws.Range("F4").Value = 2846
$ws.Range("F5").Value = 1753
$ws.Range("F7").Value = 7
$ws.Range("F8").Value = 763
$ws.Range("F11").Value = 444
$ws.Range("F13").Value = 307
$ws.Range("F14").Value = 13
$ws.Range("F16").Value = 82
$ws.Range("F17").Value = 1056
$ws.Range("F19").Value = 691
$ws.Range("F20").Value = 765
$ws.Range("F22").Value = 111
$ws.Range("F23").Value = 111
$ws.Range("F25").Value = 17
$ws.Range("F28").Value = 586
$ws.Range("F29").Value = 72
$ws.Range("F30").Value = 1670
$ws.Range("F31").Value = 1747
$ws.Range("F32").Value = 431
$ws.Range("F34").Value = 2440
$ws.Range("F35").Value = 429
$ws.Range("F38").Value = 53
$ws.Range("F39").Value = 637
$ws.Range("F41").Value = 75
$ws.Range("F44").Value = 1536
$ws.Range("F45").Value = 244
$ws.Range("F49").Value = 126
